# Docx writer: Use Compact style for Plain block elements.
# This differentiates between tight and loose lists.
#
# 1. Give the "Normal" style explicit paragraph spacing (before/after = 9pt,
#    i.e. 180 twips) so it keeps acting as the "loose" baseline.
# 2. Add a new "Compact" paragraph style, based on "Normal", with much
#    tighter spacing (before/after = 1.8pt, i.e. 36 twips) for use by tight
#    list items / Plain block elements.

$d = $word.ActiveDocument

$wdStyleTypeParagraph = 1

$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceBefore = 9
$normal.ParagraphFormat.SpaceAfter = 9

$compact = $d.Styles.Add("Compact", $wdStyleTypeParagraph)
$compact.BaseStyle = "Normal"
$compact.ParagraphFormat.SpaceBefore = 1.8
$compact.ParagraphFormat.SpaceAfter = 1.8
$compact.QuickStyle = $true
